# Bug fix: column D (过期时间 / expiry date) held a mix of real date
# numbers (e.g. row 3/5) and plain text strings that merely looked like
# dates (e.g. "2018年02月14日" in row 4/6). Normalize the whole column to
# real date-time values using a single, consistent number format so every
# cell in D3:D6 is an actual Excel date rather than some dates / some text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the unified date/time format to the whole range first so the
# subsequent numeric writes are interpreted (and stored) as dates rather
# than being coerced back to text by a pre-existing "store as text" format.
$ws.Range("D3:D6").NumberFormat = "yyyy.m.d h:mm"

# Write proper date-time serial values (replacing the old plain date / text
# values) for every row in the 过期时间 column.
$ws.Range("D3").Value2 = 41711.382638888892
$ws.Range("D4").Value2 = 41711.757638888892
$ws.Range("D5").Value2 = 41407.757638888892
$ws.Range("D6").Value2 = 41711.757638888892
